$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the worker/period data table (rows 16-64): grouped by worker,
# periods descending 2310..2304, each worker getting all 7 periods.
# Period 2310 carries the reduced value (35574); all other periods keep 46400.
$workers = @(
    @('1042419895', 'LILIAN MILENA HERRERA FIGUEROA'),
    @('8828922', 'ANTONINO RAMOS SAAVEDRA'),
    @('1050556709', 'DEIVER STIVEN RAMOS BARBA'),
    @('1002361731', 'DIDIER IVAN RAMOS BARBA'),
    @('1050554318', 'ABEL ANTONIO PADILLA ESCOBAR'),
    @('13620423', 'GELSON CUELLAR NOGUERA'),
    @('91323747', 'JORGE MEJIA BELTRAN'),
)
$periods = @('2310', '2309', '2308', '2307', '2306', '2305', '2304')

$row = 16
foreach ($worker in $workers) {
    $doc = $worker[0]
    $name = $worker[1]
    foreach ($period in $periods) {
        if ($period -eq "2310") { $valorMora = 35574 } else { $valorMora = 46400 }
        $ws.Cells.Item($row, 3).Value = $doc
        $ws.Cells.Item($row, 4).Value = $name
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = $valorMora
        $row++
    }
}
